# WGIC_rating_manual.xlsx - apply the edits captured by the target diff.
#
# Context (commit message): a topological sorter / DAG was added to drive
# rating order, plus multiprocessing to rate several rating steps
# concurrently. As part of that work two rating-factor sheets were renamed
# to match the new (shorter) naming convention used by the rating code, and
# the workbook was left with a different sheet active/selected than before.

$wb = $excel.ActiveWorkbook

# --- Rename the two rating-factor sheets -----------------------------------
# sheet4: protectionclass_constructiontyp -> protectclass_constr_factor
# sheet8: five_year_claim_free_discount_f -> five_year_claim_free_factor
$wb.Worksheets.Item(4).Name = "protectclass_constr_factor"
$wb.Worksheets.Item(8).Name = "five_year_claim_free_factor"

# --- Move the active tab / selection ---------------------------------------
# Previously "multi_policy_discount_factor" (9th tab) was the active sheet
# (selection M15). Now "protectclass_constr_factor" (4th tab, just renamed
# above) is the active sheet, with the selection moved to D26.
$ws4 = $wb.Worksheets.Item(4)
$ws4.Activate() | Out-Null
$ws4.Range("D26").Select() | Out-Null
